$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of columns C (group-name) and D (group-code) for every
# used row, including the header row. This matches the shared-string
# reordering seen in the diff (group-code now precedes group-name).
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value()
    $dVal = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 3).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $cVal
}
